$d = $word.ActiveDocument

# 1. Insert a new paragraph "BIN PHishCODE:Jim1975" right before the "Gmail: Jamilchik1975" paragraph.
$gmail = $d.Content
$gmail.Find.Execute("Gmail: Jamilchik1975", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
if ($gmail.Find.Found) {
    $gmailPara = $gmail.Paragraphs(1)
    $insertionPoint = $gmailPara.Range.Start
    $ins = $d.Range($insertionPoint, $insertionPoint)
    $ins.InsertBefore("BIN PHishCODE:Jim1975`r")
}

# 2. Normalize the "GAKhome: 04522966" paragraph into a single clean run,
#    dropping the stray proofErr spell-check markers that split it in two runs.
$gak = $d.Content
$gak.Find.Execute("GAKhome: 04522966", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
if ($gak.Find.Found) {
    $gakPara = $gak.Paragraphs(1)
    $gakRange = $gakPara.Range
    $insertionPoint2 = $gakRange.Start
    $gakRange.Delete()
    $ins2 = $d.Range($insertionPoint2, $insertionPoint2)
    $ins2.InsertBefore("GAKhome: 04522966`r")
}
